$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure the Price/Volume columns keep their original text format so that
# numeric-looking strings (e.g. "300.97") are not auto-converted to numbers.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.805.38'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.363.94'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.99%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '300.97'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '95.76'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.490'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.91%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '33.91'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -1.71%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0786'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +2.57%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '18.27'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -3.46%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.74'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.738.65'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +2.19%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.364.73'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.798'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.61%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '42.793.59'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.12'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.28'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +1.99%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0₃0885'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -0.90%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '67.92'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '234.81'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -1.73%  '
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '24.73'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +1.59%  '
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.18'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '31.44'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -2.66%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.03'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0733'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +4.82%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '17.17'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -4.25%  '
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +4.92%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.85'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +4.02%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -2.61%  '
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -1.40%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.78'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +1.69%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '22.34'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +6.73%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '117.22'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -29.45%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.935.21'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.28%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0279'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +1.87%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.23%  '
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -9.19%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.595.30'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +1.87%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +1.71%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '71.95'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '51.87'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -2.97%  '
